$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1333
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H5").Value = 356.2
$ws.Range("I5").Value = 344.6
$ws.Range("J5").Value = 367.8
$ws.Range("K5").Value = 344.6
$ws.Range("L5").Value = 367.8
$ws.Range("M5").Value = -229.6
$ws.Range("N5").Value = -597.8
$ws.Range("H6").Value = 203.75
$ws.Range("I6").Value = 203.75
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 611.25
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -499.25
$ws.Range("N6").ClearContents()
$ws.Range("H18").Value = 2965.158
$ws.Range("I18").Value = 2450
$ws.Range("K18").Value = 2450
$ws.Range("M18").Value = -2166
$ws.Range("H40").Value = 4487.2896
$ws.Range("J40").Value = 4784.88
$ws.Range("L40").Value = 4784.88
$ws.Range("N40").Value = -5134.88
$ws.Range("H64").Value = 7932.8667
$ws.Range("J64").Value = 8073.4814
$ws.Range("L64").Value = 8073.4814
$ws.Range("N64").Value = -8569.481400000001
$ws.Range("H67").Value = 7932.8667
$ws.Range("J67").Value = 8073.4814
$ws.Range("L67").Value = 8073.4814
$ws.Range("N67").Value = -9789.481400000001
$ws.Range("H70").Value = 6974.925
$ws.Range("J70").Value = 6974.925
$ws.Range("L70").Value = 20924.775
$ws.Range("N70").Value = -21464.775
$ws.Range("H73").Value = 6974.925
$ws.Range("J73").Value = 6974.925
$ws.Range("L73").Value = 20924.775
$ws.Range("N73").Value = -22796.775
$ws.Range("H74").Value = 7388.7856
$ws.Range("I74").Value = 4726.5
$ws.Range("K74").Value = 4726.5
$ws.Range("M74").Value = -3790.5
$ws.Range("H77").Value = 7388.7856
$ws.Range("I77").Value = 4726.5
$ws.Range("K77").Value = 23632.5
$ws.Range("M77").Value = -18952.5
$ws.Range("H80").Value = 5178
$ws.Range("I80").Value = 986.6923
$ws.Range("J80").Value = 8383.117
$ws.Range("K80").Value = 2960.0769
$ws.Range("L80").Value = 25149.351
$ws.Range("M80").Value = -1962.0769
$ws.Range("N80").Value = -27145.351
$ws.Range("H83").Value = 5178
$ws.Range("I83").Value = 986.6923
$ws.Range("J83").Value = 8383.117
$ws.Range("K83").Value = 8880.2307
$ws.Range("L83").Value = 75448.053
$ws.Range("M83").Value = -3888.2307
$ws.Range("N83").Value = -85432.053
$ws.Range("H98").Value = 1974.238
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1974.238
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 18521054
$ws.Range("I132").Value = 18870478
$ws.Range("K132").Value = 56611434
$ws.Range("M132").Value = -56608904
$ws.Range("H137").Value = 42421.76
$ws.Range("I137").Value = 108487.94
$ws.Range("J137").Value = 3693.3103
$ws.Range("K137").Value = 325463.82
$ws.Range("L137").Value = 11079.9309
$ws.Range("M137").Value = -322913.82
$ws.Range("N137").Value = -16179.9309
$ws.Range("H138").Value = 2992.3232
$ws.Range("I138").Value = 804.375
$ws.Range("J138").Value = 3414.0964
$ws.Range("K138").Value = 2413.125
$ws.Range("L138").Value = 10242.2892
$ws.Range("M138").Value = 2726.875
$ws.Range("N138").Value = -20522.2892
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360
$ws.Range("H141").Value = 4917.5864
$ws.Range("I141").Value = 5662.909
$ws.Range("J141").Value = 2575.1428
$ws.Range("K141").Value = 16988.727
$ws.Range("L141").Value = 7725.428400000001
$ws.Range("M141").Value = -11808.727
$ws.Range("N141").Value = -18085.4284

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 164.81818
$ws.Range("I4").Value = 181.1
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 181.1
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = -65.09999999999999
$ws.Range("N4").Value = -234
$ws.Range("H5").Value = 404.33334
$ws.Range("I5").Value = 393.5
$ws.Range("J5").Value = 426
$ws.Range("K5").Value = 393.5
$ws.Range("L5").Value = 426
$ws.Range("M5").Value = -281.5
$ws.Range("N5").Value = -650
$ws.Range("H32").Value = 4114.4067
$ws.Range("I32").Value = 2473.3247
$ws.Range("K32").Value = 2473.3247
$ws.Range("M32").Value = -2186.3247
$ws.Range("H53").Value = 16867.8
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 16867.8
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 16867.8
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -18231.8
$ws.Range("H61").Value = 3029.96
$ws.Range("I61").Value = 2939.8
$ws.Range("K61").Value = 2939.8
$ws.Range("M61").Value = -2727.8
$ws.Range("H122").Value = 1491196.9
$ws.Range("I122").Value = 3045.4
$ws.Range("J122").Value = 5211576
$ws.Range("K122").Value = 9136.200000000001
$ws.Range("L122").Value = 15634728
$ws.Range("M122").Value = -6686.200000000001
$ws.Range("N122").Value = -15639628
$ws.Range("H132").Value = 1648.1094
$ws.Range("I132").Value = 1352.6444
$ws.Range("K132").Value = 4057.933199999999
$ws.Range("M132").Value = -1527.933199999999
$ws.Range("H136").Value = 3029.96
$ws.Range("I136").Value = 2939.8
$ws.Range("K136").Value = 8819.400000000001
$ws.Range("M136").Value = -6269.400000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 404.33334
$ws.Range("I4").Value = 393.5
$ws.Range("J4").Value = 426
$ws.Range("K4").Value = 393.5
$ws.Range("L4").Value = 426
$ws.Range("M4").Value = -278.5
$ws.Range("N4").Value = -656
$ws.Range("H86").Value = 11124444
$ws.Range("I86").Value = 20020998
$ws.Range("J86").Value = 3750
$ws.Range("K86").Value = 20020998
$ws.Range("L86").Value = 3750
$ws.Range("M86").Value = -20019875
$ws.Range("N86").Value = -5996
$ws.Range("H89").Value = 11124444
$ws.Range("I89").Value = 20020998
$ws.Range("J89").Value = 3750
$ws.Range("K89").Value = 100104990
$ws.Range("L89").Value = 18750
$ws.Range("M89").Value = -100099374
$ws.Range("N89").Value = -29982
$ws.Range("H94").Value = 3578327.2
$ws.Range("I94").Value = 5000320
$ws.Range("K94").Value = 5000320
$ws.Range("M94").Value = -4999869
$ws.Range("H105").Value = 3679743.5
$ws.Range("I105").Value = 3909540
$ws.Range("K105").Value = 3909540
$ws.Range("M105").Value = -3907793
$ws.Range("H122").Value = 75170
$ws.Range("J122").Value = 75170
$ws.Range("L122").Value = 75170
$ws.Range("N122").Value = -84970
$ws.Range("H132").Value = 85085
$ws.Range("J132").Value = 85085
$ws.Range("L132").Value = 85085
$ws.Range("N132").Value = -95205

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 324.52
$ws.Range("J7").Value = 526.6667
$ws.Range("L7").Value = 526.6667
$ws.Range("N7").Value = -752.6667
$ws.Range("H22").Value = 576.4286
$ws.Range("I22").Value = 347.2
$ws.Range("K22").Value = 347.2
$ws.Range("M22").Value = 2.800000000000011
$ws.Range("H31").Value = 16870.307
$ws.Range("I31").Value = 1428.7693
$ws.Range("J31").Value = 28022.527
$ws.Range("K31").Value = 1428.7693
$ws.Range("L31").Value = 28022.527
$ws.Range("M31").Value = -1133.7693
$ws.Range("N31").Value = -28612.527
$ws.Range("H34").Value = 16870.307
$ws.Range("I34").Value = 1428.7693
$ws.Range("J34").Value = 28022.527
$ws.Range("K34").Value = 1428.7693
$ws.Range("L34").Value = 28022.527
$ws.Range("M34").Value = -1226.7693
$ws.Range("N34").Value = -28426.527
$ws.Range("H134").Value = 2627.7886
$ws.Range("I134").Value = 2157.739
$ws.Range("K134").Value = 6473.217000000001
$ws.Range("M134").Value = -3938.217000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 264.7353
$ws.Range("I2").Value = 96.13333
$ws.Range("J2").Value = 397.8421
$ws.Range("K2").Value = 576.79998
$ws.Range("L2").Value = 2387.0526
$ws.Range("M2").Value = -463.79998
$ws.Range("N2").Value = -2613.0526
$ws.Range("H4").Value = 5568659
$ws.Range("I4").Value = 6124036
$ws.Range("J4").Value = 200013.33
$ws.Range("K4").Value = 18372108
$ws.Range("L4").Value = 600039.99
$ws.Range("M4").Value = -18371996
$ws.Range("N4").Value = -600263.99
$ws.Range("H37").Value = 53249.9
$ws.Range("J37").Value = 53249.9
$ws.Range("L37").Value = 159749.7
$ws.Range("N37").Value = -159973.7
$ws.Range("H80").Value = 8492.5
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 6985
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 20955
$ws.Range("M80").Value = -29064
$ws.Range("N80").Value = -22827
$ws.Range("H81").Value = 6832.8
$ws.Range("J81").Value = 8306.666999999999
$ws.Range("L81").Value = 24920.001
$ws.Range("N81").Value = -27166.001
$ws.Range("H83").Value = 8492.5
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 6985
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 62865
$ws.Range("M83").Value = -85320
$ws.Range("N83").Value = -72225
$ws.Range("H84").Value = 6832.8
$ws.Range("J84").Value = 8306.666999999999
$ws.Range("L84").Value = 74760.003
$ws.Range("N84").Value = -85992.003
$ws.Range("H107").Value = 417.16666
$ws.Range("I107").Value = 223.625
$ws.Range("J107").Value = 572
$ws.Range("K107").Value = 670.875
$ws.Range("L107").Value = 1716
$ws.Range("M107").Value = 1249.125
$ws.Range("N107").Value = -5556
$ws.Range("H129").Value = 1071.3889
$ws.Range("I129").Value = 967.4545000000001
$ws.Range("J129").Value = 1234.7142
$ws.Range("K129").Value = 2902.3635
$ws.Range("L129").Value = 3704.1426
$ws.Range("M129").Value = 2097.6365
$ws.Range("N129").Value = -13704.1426
$ws.Range("H132").Value = 2019.3823
$ws.Range("J132").Value = 2485.9167
$ws.Range("L132").Value = 22373.2503
$ws.Range("N132").Value = -27433.2503
$ws.Range("H139").Value = 607.7
$ws.Range("I139").Value = 492.3684
$ws.Range("J139").Value = 2799
$ws.Range("K139").Value = 1477.1052
$ws.Range("L139").Value = 8397
$ws.Range("M139").Value = 3662.8948
$ws.Range("N139").Value = -18677

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 11069.9
$ws.Range("I2").Value = 1450.5714
$ws.Range("J2").Value = 33515
$ws.Range("K2").Value = 1450.5714
$ws.Range("L2").Value = 33515
$ws.Range("M2").Value = -1337.5714
$ws.Range("N2").Value = -33741
$ws.Range("H4").Value = 13075.75
$ws.Range("I4").Value = 17267.666
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 17267.666
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -17155.666
$ws.Range("N4").Value = -724
$ws.Range("H109").Value = 46357.8
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 46357.8
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 46357.8
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -48437.8
$ws.Range("H122").Value = 411445.25
$ws.Range("I122").Value = 456338.16
$ws.Range("K122").Value = 1369014.48
$ws.Range("M122").Value = -1366564.48

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9700.4
$ws.Range("I7").Value = 8375.75
$ws.Range("K7").Value = 8375.75
$ws.Range("M7").Value = -8263.75
$ws.Range("H18").Value = 15000
$ws.Range("J18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15344
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H40").Value = 14999
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 8666.666999999999
$ws.Range("J46").Value = 9800
$ws.Range("L46").Value = 9800
$ws.Range("N46").Value = -10176
$ws.Range("H119").Value = 66128.336
$ws.Range("J119").Value = 66128.336
$ws.Range("L119").Value = 66128.336
$ws.Range("N119").Value = -75804.336
$ws.Range("H126").Value = 9700.4
$ws.Range("I126").Value = 8375.75
$ws.Range("K126").Value = 25127.25
$ws.Range("M126").Value = -22657.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 7500
$ws.Range("J19").Value = 7500
$ws.Range("L19").Value = 7500
$ws.Range("N19").Value = -7848
$ws.Range("H62").Value = 8581.481
$ws.Range("I62").Value = 4249.75
$ws.Range("K62").Value = 4249.75
$ws.Range("M62").Value = -3625.75
$ws.Range("H65").Value = 8581.481
$ws.Range("I65").Value = 4249.75
$ws.Range("K65").Value = 21248.75
$ws.Range("M65").Value = -18128.75
$ws.Range("H81").Value = 83334340
$ws.Range("I81").Value = 83334340
$ws.Range("K81").Value = 166668680
$ws.Range("M81").Value = -166667619
$ws.Range("H84").Value = 83334340
$ws.Range("I84").Value = 83334340
$ws.Range("K84").Value = 833343400
$ws.Range("M84").Value = -833338096
